$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Jira")

# Sheet1: add descriptions for row 5 (function prototype) and row 4 (struct)
# Note: written in this order so new shared strings are appended as
# "declare isValidInput()" then "make truck and shipement struct"
$ws1.Range("C5").Value = "declare isValidInput()"
$ws1.Range("D5").Value = "Fail"
$ws1.Range("E5").Value = "Fail"
$ws1.Range("F5").Value = "Fail"

$ws1.Range("C4").Value = "make truck and shipement struct"
$ws1.Range("D4").Value = "Fail"
$ws1.Range("E4").Value = "Fail"
$ws1.Range("F4").Value = "Fail"

# Column C on Sheet1 grows to fit the new, longer description text
$ws1.Columns.Item(3).EntireColumn.AutoFit()

# Jira sheet: widen column I to fit its content and move the selection
$ws2.Columns.Item(9).EntireColumn.AutoFit()
$ws2.Range("I10").Select()

# Update the remembered cell selection on Sheet1 last, so Sheet1 stays
# the active/selected tab, matching the source workbook
$ws1.Range("L13").Select()
